# Testing and bug fixing: drop the "NT" (not tested) markers that were
# placeholders for rows/columns not yet exercised, and update the sheet's
# saved view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B6 held "NT - not tested"; E12:H12 held "NT" — clear them back out now
# that those cases have been tested.
$ws.Range("B6").ClearContents()
$ws.Range("E12:H12").ClearContents()

# Scroll the view back to the top and leave the selection where editing
# left off.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("Q24").Select()
